# Natmi following Dr Hou advice
# Update the LR-pairs results table (rows 2-13 revised, rows 14-17 added)
# so that it contains the full 4x4 cross product of sending/target clusters
# (ECs, FAPs, M2, sCs) for the Col1a1-Ddr1 ligand-receptor pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one data row: columns A..T
$data = @(
    @("ECs", "Col1a1", "Ddr1", "ECs", 3, 1, 9.108069666666667, 27.324209, 0.00155006418458712, 0.00155006418458712, 2, 0.6666666666666666, 0.353194, 1.059582, 0.01734490358674178, 0.01734490358674178, 3.216915557848667, 28.952240020638, [double]"2.68857138349251e-05", [double]"2.68857138349251e-05"),
    @("ECs", "Col1a1", "Ddr1", "FAPs", 3, 1, 9.108069666666667, 27.324209, 0.00155006418458712, 0.00155006418458712, 3, 1, 1.890286333333333, 5.670859, 0.09282953335278148, 0.09282953335278148, 17.21685961394789, 154.951736525531, 0.000143891734922082, 0.000143891734922082),
    @("ECs", "Col1a1", "Ddr1", "M2", 3, 1, 9.108069666666667, 27.324209, 0.00155006418458712, 0.00155006418458712, 3, 1, 0.2252136666666666, 0.6756409999999999, 0.01105995383486111, 0.01105995383486111, 2.051261765885444, 18.461355892969, [double]"1.714363832260517e-05", [double]"1.714363832260517e-05"),
    @("ECs", "Col1a1", "Ddr1", "sCs", 3, 1, 9.108069666666667, 27.324209, 0.00155006418458712, 0.00155006418458712, 3, 1, 17.89429033333333, 53.68287100000001, 0.8787656092256156, 0.8787656092256156, 162.9824429915599, 1466.841986924039, 0.001362143097507507, 0.001362143097507507),
    @("FAPs", "Col1a1", "Ddr1", "ECs", 3, 1, 5771.873535333333, 17315.620606, 0.9822909543423312, 0.9822909543423313, 2, 0.6666666666666666, 0.353194, 1.059582, 0.01734490358674178, 0.01734490358674178, 2038.591101438521, 18347.31991294669, 0.01703774189719631, 0.01703774189719631),
    @("FAPs", "Col1a1", "Ddr1", "FAPs", 3, 1, 5771.873535333333, 17315.620606, 0.9822909543423312, 0.9822909543423313, 3, 1, 1.890286333333333, 5.670859, 0.09282953335278148, 0.09282953335278148, 10910.49366156895, 98194.44295412056, 0.09118561090825697, 0.09118561090825698),
    @("FAPs", "Col1a1", "Ddr1", "M2", 3, 1, 5771.873535333333, 17315.620606, 0.9822909543423312, 0.9822909543423313, 3, 1, 0.2252136666666666, 0.6756409999999999, 0.01105995383486111, 0.01105995383486111, 1299.904802428716, 11699.14322185845, 0.01086409260742784, 0.01086409260742784),
    @("FAPs", "Col1a1", "Ddr1", "sCs", 3, 1, 5771.873535333333, 17315.620606, 0.9822909543423312, 0.9822909543423313, 3, 1, 17.89429033333333, 53.68287100000001, 0.8787656092256156, 0.8787656092256156, 103283.5808085378, 929552.2272768399, 0.8632035089294501, 0.8632035089294501),
    @("M2", "Col1a1", "Ddr1", "ECs", 3, 1, 1.272029666666667, 3.816089, 0.0002164813950916887, 0.0002164813950916887, 2, 0.6666666666666666, 0.353194, 1.059582, 0.01734490358674178, 0.01734490358674178, 0.4492732460886666, 4.043459214798, [double]"3.754848926188695e-06", [double]"3.754848926188696e-06"),
    @("M2", "Col1a1", "Ddr1", "FAPs", 3, 1, 1.272029666666667, 3.816089, 0.0002164813950916887, 0.0002164813950916887, 3, 1, 1.890286333333333, 5.670859, 0.09282953335278148, 0.09282953335278148, 2.404500294494555, 21.640502650451, [double]"2.009586688592058e-05", [double]"2.009586688592058e-05"),
    @("M2", "Col1a1", "Ddr1", "M2", 3, 1, 1.272029666666667, 3.816089, 0.0002164813950916887, 0.0002164813950916887, 3, 1, 0.2252136666666666, 0.6756409999999999, 0.01105995383486111, 0.01105995383486111, 0.2864784653387777, 2.578306188048999, [double]"2.394274235820405e-06", [double]"2.394274235820405e-06"),
    @("M2", "Col1a1", "Ddr1", "sCs", 3, 1, 1.272029666666667, 3.816089, 0.0002164813950916887, 0.0002164813950916887, 3, 1, 17.89429033333333, 53.68287100000001, 0.8787656092256156, 0.8787656092256156, 22.76206816794656, 204.858613511519, 0.000190236405043759, 0.000190236405043759),
    @("sCs", "Col1a1", "Ddr1", "ECs", 3, 1, 93.67702500000001, 281.031075, 0.01594250007799006, 0.01594250007799006, 2, 0.6666666666666666, 0.353194, 1.059582, 0.01734490358674178, 0.01734490358674178, 33.08616316785, 297.7754685106501, 0.0002765211267843609, 0.0002765211267843609),
    @("sCs", "Col1a1", "Ddr1", "FAPs", 3, 1, 93.67702500000001, 281.031075, 0.01594250007799006, 0.01594250007799006, 3, 1, 1.890286333333333, 5.670859, 0.09282953335278148, 0.09282953335278148, 177.076400104825, 1593.687600943425, 0.0014799348427165, 0.0014799348427165),
    @("sCs", "Col1a1", "Ddr1", "M2", 3, 1, 93.67702500000001, 281.031075, 0.01594250007799006, 0.01594250007799006, 3, 1, 0.2252136666666666, 0.6756409999999999, 0.01105995383486111, 0.01105995383486111, 21.097346282675, 189.876116544075, 0.0001763233148748396, 0.0001763233148748396),
    @("sCs", "Col1a1", "Ddr1", "sCs", 3, 1, 93.67702500000001, 281.031075, 0.01594250007799006, 0.01594250007799006, 3, 1, 17.89429033333333, 53.68287100000001, 0.8787656092256156, 0.8787656092256156, 1676.283882912925, 15086.55494621633, 0.01400972079361436, 0.01400972079361436)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $row[$j]
    }
}
